# Updates the cryptos price-table "Price" (D) and "Volume(1h)" (E) columns
# to the latest scraped snapshot. Values are plain display strings (not
# numbers), matching the source data which uses "." as a thousands
# separator (e.g. "41.522.93") alongside plain decimals (e.g. "314.69")
# and padded percentage strings (e.g. "  +0.13%  ").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. ForceText marks values that Excel would otherwise
# auto-convert to a number (single "." decimal) so we flip the cell to Text
# format just long enough to assign the literal string, then restore General
# (matching the unstyled cells in the original workbook).
$updates = @(
    @{ Cell = 'D2'; Value = '41.522.93'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  +0.13%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '2.466.95'; ForceText = $false }
    @{ Cell = 'E4'; Value = '  -0.11%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '314.69'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +0.68%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '91.86'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -2.70%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '0.546'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  -0.40%  '; ForceText = $false }
    @{ Cell = 'E8'; Value = '  -0.20%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '0.512'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  +2.60%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '32.39'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -3.41%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.0790'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  +0.90%  '; ForceText = $false }
    @{ Cell = 'E12'; Value = '  +1.22%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '2.850.06'; ForceText = $false }
    @{ Cell = 'E13'; Value = '  -0.28%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '6.84'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -2.06%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '16.08'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +5.14%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '2.532.45'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  +1.92%  '; ForceText = $false }
    @{ Cell = 'E17'; Value = '  -2.93%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '41.498.40'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  +0.23%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '6.49'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  +2.72%  '; ForceText = $false }
    @{ Cell = 'E20'; Value = '  +2.34%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '71.58'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  +4.10%  '; ForceText = $false }
    @{ Cell = 'E22'; Value = '  -1.72%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '236.03'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -0.58%  '; ForceText = $false }
    @{ Cell = 'E24'; Value = '  -1.23%  '; ForceText = $false }
    @{ Cell = 'E25'; Value = '  -0.09%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '1.88'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -1.04%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '24.77'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  +2.85%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '2.22'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -1.33%  '; ForceText = $false }
    @{ Cell = 'E30'; Value = '  -2.82%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '156.15'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  +2.73%  '; ForceText = $false }
    @{ Cell = 'E32'; Value = '  -0.45%  '; ForceText = $false }
    @{ Cell = 'E33'; Value = '  -0.06%  '; ForceText = $false }
    @{ Cell = 'E34'; Value = '  +1.41%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '17.23'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -1.54%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '2.33'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -9.03%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '2.86'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -6.99%  '; ForceText = $false }
    @{ Cell = 'E38'; Value = '  +1.60%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '0.113'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  -0.70%  '; ForceText = $false }
    @{ Cell = 'E40'; Value = '  -4.63%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '4.04'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -4.77%  '; ForceText = $false }
    @{ Cell = 'E42'; Value = '  -0.30%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '1.955.58'; ForceText = $false }
    @{ Cell = 'E43'; Value = '  -1.60%  '; ForceText = $false }
    @{ Cell = 'E44'; Value = '  -0.82%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '18.56'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -4.50%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '2.92'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -2.39%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '9.07'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  +3.85%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '2.710.88'; ForceText = $false }
    @{ Cell = 'E48'; Value = '  -0.43%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '97.12'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -0.04%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '66.98'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -3.67%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '71.79'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -3.68%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.NumberFormat = "General"
    } else {
        $rng.Value = $u.Value
    }
}

